$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(5, 8).Value = 630.875
$ws.Cells.Item(5, 9).Value = 630.875
$ws.Cells.Item(5, 11).Value = 630.875
$ws.Cells.Item(5, 13).Value = -515.875
$ws.Cells.Item(96, 8).Value = 443
$ws.Cells.Item(96, 9).Value = 401.45456
$ws.Cells.Item(96, 11).Value = 1204.36368
$ws.Cells.Item(96, 13).Value = 168.6363200000001
$ws.Cells.Item(135, 8).Value = 1662.2593
$ws.Cells.Item(135, 9).Value = 1627.3636
$ws.Cells.Item(135, 11).Value = 14646.2724
$ws.Cells.Item(135, 13).Value = -12111.2724
$ws.Cells.Item(138, 8).Value = 3427.9434
$ws.Cells.Item(138, 9).Value = 2462.0454
$ws.Cells.Item(138, 10).Value = 4113.4194
$ws.Cells.Item(138, 11).Value = 7386.1362
$ws.Cells.Item(138, 12).Value = 12340.2582
$ws.Cells.Item(138, 13).Value = -2246.1362
$ws.Cells.Item(138, 14).Value = -22620.2582

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 5212.25
$ws.Cells.Item(32, 9).Value = 3367.1838
$ws.Cells.Item(32, 10).Value = 17560
$ws.Cells.Item(32, 11).Value = 3367.1838
$ws.Cells.Item(32, 12).Value = 17560
$ws.Cells.Item(32, 13).Value = -3080.1838
$ws.Cells.Item(32, 14).Value = -18134
$ws.Cells.Item(61, 8).Value = 12057.308
$ws.Cells.Item(61, 9).Value = 5500
$ws.Cells.Item(61, 10).Value = 14971.667
$ws.Cells.Item(61, 11).Value = 5500
$ws.Cells.Item(61, 12).Value = 14971.667
$ws.Cells.Item(61, 13).Value = -5288
$ws.Cells.Item(61, 14).Value = -15395.667
$ws.Cells.Item(97, 8).Value = 525.75
$ws.Cells.Item(97, 9).Value = 478.2353
$ws.Cells.Item(97, 11).Value = 478.2353
$ws.Cells.Item(97, 13).Value = 17.7647
$ws.Cells.Item(132, 8).Value = 732230.25
$ws.Cells.Item(132, 9).Value = 1111694.1
$ws.Cells.Item(132, 10).Value = 7799.091
$ws.Cells.Item(132, 11).Value = 3335082.3
$ws.Cells.Item(132, 12).Value = 23397.273
$ws.Cells.Item(132, 13).Value = -3332552.3
$ws.Cells.Item(132, 14).Value = -28457.273
$ws.Cells.Item(136, 8).Value = 12057.308
$ws.Cells.Item(136, 9).Value = 5500
$ws.Cells.Item(136, 10).Value = 14971.667
$ws.Cells.Item(136, 11).Value = 16500
$ws.Cells.Item(136, 12).Value = 44915.001
$ws.Cells.Item(136, 13).Value = -13950
$ws.Cells.Item(136, 14).Value = -50015.001

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(35, 8).Value = 93997.75
$ws.Cells.Item(35, 10).Value = 93997.75
$ws.Cells.Item(35, 12).Value = 93997.75
$ws.Cells.Item(35, 14).Value = -94617.75
$ws.Cells.Item(42, 8).Value = 0
$ws.Cells.Item(42, 10).Value = 0
$ws.Cells.Item(42, 12).Value = 0
$ws.Cells.Item(42, 14).ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 7808.3076
$ws.Cells.Item(31, 9).Value = 1487.6875
$ws.Cells.Item(31, 10).Value = 10617.473
$ws.Cells.Item(31, 11).Value = 1487.6875
$ws.Cells.Item(31, 12).Value = 10617.473
$ws.Cells.Item(31, 13).Value = -1192.6875
$ws.Cells.Item(31, 14).Value = -11207.473
$ws.Cells.Item(34, 8).Value = 7808.3076
$ws.Cells.Item(34, 9).Value = 1487.6875
$ws.Cells.Item(34, 10).Value = 10617.473
$ws.Cells.Item(34, 11).Value = 1487.6875
$ws.Cells.Item(34, 12).Value = 10617.473
$ws.Cells.Item(34, 13).Value = -1285.6875
$ws.Cells.Item(34, 14).Value = -11021.473
$ws.Cells.Item(50, 8).Value = 119994.2
$ws.Cells.Item(50, 10).Value = 119994.2
$ws.Cells.Item(50, 12).Value = 119994.2
$ws.Cells.Item(50, 14).Value = -121244.2
$ws.Cells.Item(51, 8).Value = 102522.25
$ws.Cells.Item(51, 9).Value = 20000
$ws.Cells.Item(51, 10).Value = 130029.664
$ws.Cells.Item(51, 11).Value = 20000
$ws.Cells.Item(51, 12).Value = 130029.664
$ws.Cells.Item(51, 14).Value = -131501.664
$ws.Cells.Item(51, 13).Value = -19264
$ws.Cells.Item(58, 8).Value = 1126123.1
$ws.Cells.Item(58, 9).Value = 1544794.9
$ws.Cells.Item(58, 11).Value = 1544794.9
$ws.Cells.Item(58, 13).Value = -1544591.9
$ws.Cells.Item(59, 8).Value = 40114.5
$ws.Cells.Item(59, 10).Value = 40114.5
$ws.Cells.Item(59, 12).Value = 40114.5
$ws.Cells.Item(59, 14).Value = -42404.5
$ws.Cells.Item(60, 8).Value = 75254.45
$ws.Cells.Item(60, 10).Value = 78829.89999999999
$ws.Cells.Item(60, 12).Value = 78829.89999999999
$ws.Cells.Item(60, 14).Value = -79851.89999999999
$ws.Cells.Item(61, 8).Value = 102522.25
$ws.Cells.Item(61, 9).Value = 20000
$ws.Cells.Item(61, 10).Value = 130029.664
$ws.Cells.Item(61, 11).Value = 20000
$ws.Cells.Item(61, 12).Value = 130029.664
$ws.Cells.Item(61, 14).Value = -130725.664
$ws.Cells.Item(61, 13).Value = -19652
$ws.Cells.Item(136, 8).Value = 1126123.1
$ws.Cells.Item(136, 9).Value = 1544794.9
$ws.Cells.Item(136, 11).Value = 4634384.699999999
$ws.Cells.Item(136, 13).Value = -4631834.699999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(7, 8).Value = 6833554
$ws.Cells.Item(7, 9).Value = 10000062
$ws.Cells.Item(7, 10).Value = 5250300
$ws.Cells.Item(7, 11).Value = 30000186
$ws.Cells.Item(7, 12).Value = 15750900
$ws.Cells.Item(7, 13).Value = -30000074
$ws.Cells.Item(7, 14).Value = -15751124
$ws.Cells.Item(12, 8).Value = 21.181818
$ws.Cells.Item(12, 9).Value = 18.25
$ws.Cells.Item(12, 10).Value = 22.857143
$ws.Cells.Item(12, 11).Value = 54.75
$ws.Cells.Item(12, 12).Value = 68.57142899999999
$ws.Cells.Item(12, 13).Value = 118.25
$ws.Cells.Item(12, 14).Value = -414.571429
$ws.Cells.Item(109, 8).Value = 4712.095
$ws.Cells.Item(109, 10).Value = 5904.2144
$ws.Cells.Item(109, 12).Value = 17712.6432
$ws.Cells.Item(109, 14).Value = -19792.6432
$ws.Cells.Item(113, 8).Value = 3575.1667
$ws.Cells.Item(113, 9).Value = 0
$ws.Cells.Item(113, 11).Value = 0
$ws.Cells.Item(113, 13).ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(80, 8).Value = 244304.72
$ws.Cells.Item(80, 9).Value = 420699.9
$ws.Cells.Item(80, 10).Value = 9111.111000000001
$ws.Cells.Item(80, 11).Value = 420699.9
$ws.Cells.Item(80, 12).Value = 9111.111000000001
$ws.Cells.Item(80, 13).Value = -419701.9
$ws.Cells.Item(80, 14).Value = -11107.111
$ws.Cells.Item(83, 8).Value = 244304.72
$ws.Cells.Item(83, 9).Value = 420699.9
$ws.Cells.Item(83, 10).Value = 9111.111000000001
$ws.Cells.Item(83, 11).Value = 2103499.5
$ws.Cells.Item(83, 12).Value = 45555.55500000001
$ws.Cells.Item(83, 13).Value = -2098507.5
$ws.Cells.Item(83, 14).Value = -55539.55500000001
$ws.Cells.Item(97, 8).Value = 1532.36
$ws.Cells.Item(97, 9).Value = 933.4545000000001
$ws.Cells.Item(97, 11).Value = 933.4545000000001
$ws.Cells.Item(97, 13).Value = -437.4545000000001
$ws.Cells.Item(107, 8).Value = 578.2727
$ws.Cells.Item(107, 9).Value = 337.42856
$ws.Cells.Item(107, 11).Value = 337.42856
$ws.Cells.Item(107, 13).Value = 1582.57144
$ws.Cells.Item(126, 8).Value = 1193080.4
$ws.Cells.Item(126, 9).Value = 1193080.4
$ws.Cells.Item(126, 11).Value = 3579241.2
$ws.Cells.Item(126, 13).Value = -3576771.2

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 80153
$ws.Cells.Item(22, 10).Value = 3499.0833
$ws.Cells.Item(22, 12).Value = 3499.0833
$ws.Cells.Item(22, 14).Value = -4089.0833
$ws.Cells.Item(27, 8).Value = 80153
$ws.Cells.Item(27, 10).Value = 3499.0833
$ws.Cells.Item(27, 12).Value = 3499.0833
$ws.Cells.Item(27, 14).Value = -3713.0833
$ws.Cells.Item(40, 8).Value = 3438.6191
$ws.Cells.Item(40, 9).Value = 2691.889
$ws.Cells.Item(40, 10).Value = 3998.6667
$ws.Cells.Item(40, 11).Value = 2691.889
$ws.Cells.Item(40, 12).Value = 3998.6667
$ws.Cells.Item(40, 13).Value = -2555.889
$ws.Cells.Item(40, 14).Value = -4270.6667
$ws.Cells.Item(46, 8).Value = 2914.84
$ws.Cells.Item(46, 9).Value = 1263.4445
$ws.Cells.Item(46, 11).Value = 1263.4445
$ws.Cells.Item(46, 13).Value = -1075.4445
$ws.Cells.Item(100, 8).Value = 5928.269
$ws.Cells.Item(100, 9).Value = 1508.5714
$ws.Cells.Item(100, 11).Value = 1508.5714
$ws.Cells.Item(100, 13).Value = -967.5714
$ws.Cells.Item(132, 8).Value = 939196.1
$ws.Cells.Item(132, 9).Value = 1238687.8
$ws.Cells.Item(132, 10).Value = 7444.222
$ws.Cells.Item(132, 11).Value = 3716063.4
$ws.Cells.Item(132, 12).Value = 22332.666
$ws.Cells.Item(132, 13).Value = -3713533.4
$ws.Cells.Item(132, 14).Value = -27392.666
$ws.Cells.Item(136, 8).Value = 8738.308000000001
$ws.Cells.Item(136, 9).Value = 7359.8
$ws.Cells.Item(136, 11).Value = 22079.4
$ws.Cells.Item(136, 13).Value = -19529.4

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(100, 8).Value = 1849.5
$ws.Cells.Item(100, 9).Value = 1865.0526
$ws.Cells.Item(100, 11).Value = 3730.1052
$ws.Cells.Item(100, 13).Value = -3189.1052
$ws.Cells.Item(132, 8).Value = 14595106
$ws.Cells.Item(132, 9).Value = 1281624.6
$ws.Cells.Item(132, 10).Value = 44550440
$ws.Cells.Item(132, 11).Value = 3844873.8
$ws.Cells.Item(132, 12).Value = 133651320
$ws.Cells.Item(132, 13).Value = -3842343.8
$ws.Cells.Item(132, 14).Value = -133656380
$ws.Cells.Item(136, 8).Value = 12289389
$ws.Cells.Item(136, 9).Value = 13605217
$ws.Cells.Item(136, 11).Value = 40815651
$ws.Cells.Item(136, 13).Value = -40813101
